$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains weekly price records (Fruta/Hortalizas) for
# "Agrícola del Norte S.A. de Arica - Tuna". This commit refreshes the
# weekly logic: the date/quality/volume/price columns (D, L, M, N, O, P,
# Q, R, S, T) for rows 2-18 are updated to reflect the new week's values.

# Row 2: was before-row 2 <- source before-row 13
$ws.Range("D2").Value = 45028
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 1075
# Row 3: was before-row 3 <- source before-row 5
$ws.Range("D3").Value = 44993
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 130
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25462
$ws.Range("S3").Value = 1273
# Row 4: was before-row 4 <- source before-row 17
$ws.Range("D4").Value = 45014
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("S4").Value = 1225
# Row 5: was before-row 5 <- source before-row 2
$ws.Range("D5").Value = 45021
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 22000
$ws.Range("O5").Value = 23000
$ws.Range("P5").Value = 22500
$ws.Range("S5").Value = 1125
# Row 6: was before-row 6 <- source before-row 3
$ws.Range("D6").Value = 44650
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 31000
$ws.Range("O6").Value = 32000
$ws.Range("P6").Value = 31500
$ws.Range("S6").Value = 1575
# Row 7: was before-row 7 <- source before-row 4
$ws.Range("D7").Value = 44650
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("S7").Value = 1475
# Row 8: was before-row 8 <- source before-row 16
$ws.Range("D8").Value = 44979
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 250
# Row 9: was before-row 9 <- source before-row 8
$ws.Range("D9").Value = 44636
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 29000
$ws.Range("P9").Value = 29500
$ws.Range("S9").Value = 1475
# Row 10: unchanged
# Row 11: was before-row 11 <- source before-row 6
$ws.Range("D11").Value = 44679
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 29000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29500
$ws.Range("Q11").Value = "$/caja 20 kilos"
$ws.Range("S11").Value = 1475
$ws.Range("T11").Value = 20
# Row 12: was before-row 12 <- source before-row 7
$ws.Range("D12").Value = 44679
$ws.Range("L12").Value = "Tercera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 24500
$ws.Range("Q12").Value = "$/caja 20 kilos"
$ws.Range("S12").Value = 1225
$ws.Range("T12").Value = 20
# Row 13: was before-row 13 <- source before-row 11
$ws.Range("D13").Value = 44965
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 34000
$ws.Range("O13").Value = 35000
$ws.Range("P13").Value = 34600
$ws.Range("Q13").Value = "$/caja 18 kilos"
$ws.Range("S13").Value = 1922
$ws.Range("T13").Value = 18
# Row 14: was before-row 14 <- source before-row 12
$ws.Range("D14").Value = 44965
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 32000
$ws.Range("O14").Value = 33000
$ws.Range("P14").Value = 32333
$ws.Range("Q14").Value = "$/caja 18 kilos"
$ws.Range("S14").Value = 1796
$ws.Range("T14").Value = 18
# Row 15: unchanged
# Row 16: was before-row 16 <- source before-row 14
$ws.Range("D16").Value = 44671
$ws.Range("M16").Value = 200
# Row 17: was before-row 17 <- source before-row 18
$ws.Range("D17").Value = 44972
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 27000
$ws.Range("O17").Value = 28000
$ws.Range("P17").Value = 27429
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 1524
$ws.Range("T17").Value = 18
# Row 18: was before-row 18 <- source before-row 9
$ws.Range("D18").Value = 44643
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 28000
$ws.Range("O18").Value = 30000
$ws.Range("P18").Value = 29000
$ws.Range("Q18").Value = "$/caja 20 kilos"
$ws.Range("R18").Value = "Región de Coquimbo"
$ws.Range("S18").Value = 1450
$ws.Range("T18").Value = 20

Write-Output "Weekly Fruta/Hortalizas update applied."
